$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Apply formatting first (copy from existing cells that already carry
# the target style), then write values so the copied number format / text
# alignment survives the write. ----

# H1:M1 use the same right-aligned "header label" style as A16/A26 ("s=11").
$ws.Range("A16").Copy()
$ws.Range("H1:M1").PasteSpecial(-4122)   # xlPasteFormats

# I5:K5 pick up the bordered-cell style already used across row 5 (H5, "s=10").
$ws.Range("H5").Copy()
$ws.Range("I5:K5").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---- New row 1 (clue markers + the little "guess" grid on the right) ----
$ws.Range("H1").Value = "1,1"
$ws.Range("I1").Value = "?"
$ws.Range("J1").Value = "?"
$ws.Range("K1").Value = 0
$ws.Range("L1").Value = "?"
$ws.Range("M1").Value = 0

$ws.Range("T1").Value = 1
$ws.Range("U1").Value = 0
$ws.Range("V1").Value = 1
$ws.Range("W1").Value = 0
$ws.Range("X1").Value = 0
$ws.Range("Z1").Value = 20

# ---- Row 2 additions ----
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("O2").Value = 18

$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 0
$ws.Range("Z2").Value = 18

# ---- Row 3 additions ----
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 10

$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 1
$ws.Range("Z3").Value = 17

# ---- Row 4 additions ----
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 0
$ws.Range("Z4").Value = 10

# ---- Row 5 additions ----
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0

$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 1
$ws.Range("Z5").Value = 9

# ---- Row 6 additions ----
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 1
$ws.Range("Z6").Value = 5

# ---- Row 8 additions ----
$ws.Range("T8").Value = 16
$ws.Range("U8").Value = 8
$ws.Range("V8").Value = 4
$ws.Range("W8").Value = 2
$ws.Range("X8").Value = 1

# ---- View: scroll back to the top and select the new "?" total cell ----
$ws.Activate()
$ws.Range("O2").Select()
